$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, "walkingToRunning", 10.46308040618896, -2.656880855560303, -2.153079509735107, -0.5190548883925257, -1.450479206421397, -0.2456109922566745),
    @(100, "walkingToRunning", 5.822948455810547, 0.93921160697937, -1.708410024642944, -0.2103044650108719, -1.79173934116638, -0.6462500674261389),
    @(200, "walkingToRunning", 76.90617370605469, -20.47752571105957, -10.06823921203613, 2.22406245532224, -1.244166610052261, -0.3580538825594288),
    @(300, "walkingToRunning", 21.08987617492676, -20.16287040710449, 30.18997001647949, 1.525593809515442, -1.465595992349037, 0.216453205767295),
    @(400, "walkingToRunning", 10.72996425628662, -3.895856618881226, 4.558961868286133, -4.037772781986149, -3.096113548004377, 2.077658006398807),
    @(500, "walkingToRunning", 16.64472198486328, -21.42762756347656, 6.525417804718018, -7.127989895909813, 1.020252332412979, 0.7169627440061292),
    @(600, "walkingToRunning", -22.31200981140137, 9.699769020080566, -17.87223434448242, 0.1363662918695203, -0.7235090681117633, 1.348089827050427),
    @(700, "walkingToRunning", -1.527544975280761, -2.107955932617188, -5.828549861907959, 6.111373266727803, -3.888724038926818, -0.6717052065211169),
    @(800, "walkingToRunning", 25.63505172729492, -27.48458480834961, 41.41721725463867, 2.017424072293039, -1.605387079201143, -3.74285111324393),
    @(900, "walkingToRunning", 9.120732307434082, -3.210182666778564, 14.09486484527588, -3.043187412426637, -0.6237803290216224, 1.956819472553102),
    @(1000, "walkingToRunning", 10.57754039764404, 3.087275743484497, 21.89474105834961, -10.0801253147263, -3.291420891988214, 5.852694562870816),
    @(1100, "walkingToRunning", -59.92874145507812, -16.7760124206543, -37.18990707397461, -1.435964323633735, -3.116479964565011, -1.45930368780239),
    @(1200, "walkingToRunning", 18.34011650085449, 2.91422700881958, -21.99715805053711, 4.520515643435422, -6.424601422797, -1.258769486447896),
    @(1300, "walkingToRunning", 8.823348999023438, -6.199845790863037, -2.462417125701904, 3.804642815384075, -2.879516224209357, 0.2698007882070206),
    @(1400, "walkingToRunning", 36.45148849487305, -4.319328308105469, 22.25251770019531, 2.50941535251603, 3.106572741227184, 1.131389019300618),
    @(1500, "walkingToRunning", 11.48440456390381, 5.108338356018066, 7.226218223571777, -5.14737601803366, 0.1312354588679825, 5.117585698477669),
    @(1600, "walkingToRunning", 1.909540176391602, -24.70354461669922, -9.863640785217283, -2.502922200470437, 3.250862176469831, -2.965100693188241),
    @(1700, "walkingToRunning", 13.16304683685303, 1.566379547119141, -18.78874206542969, 2.697737143194075, -2.115886496125361, -4.09576786507806),
    @(1800, "walkingToRunning", -26.10807037353516, -18.4659252166748, -26.81548118591309, 5.033141990359741, -4.757362827122633, -1.592488460403557),
    @(1900, "walkingToRunning", 0.605715274810791, 0.4431395530700683, 13.34832572937012, 4.163411040957872, 2.378721077665152, -1.909676737064927),
    @(2000, "walkingToRunning", -0.1363797187805175, 12.19958686828613, 1.961617946624756, -2.440776584817322, 2.270052646561497, 2.502104572254968),
    @(2100, "walkingToRunning", 23.61079406738281, -17.57328033447266, -2.193590641021729, -4.33349856712834, 5.117880192592036, -1.15704054626636),
    @(2200, "walkingToRunning", -16.4272632598877, 6.268374919891357, -29.44120216369629, 0.3903590963899126, 12.85140922601278, -6.864528203182092),
    @(2300, "walkingToRunning", -22.77685928344727, -11.73852729797363, -5.511547565460205, 2.438172858395985, 7.542828546153372, -8.415844210617843),
    @(2400, "walkingToRunning", -20.23676681518555, -30.38065719604492, 13.72379970550537, -0.09949780518202356, -6.494732479397333, -8.307635595472599),
    @(2500, "walkingToRunning", -0.5074601173400879, 2.13600492477417, -3.543452739715576, -1.967174616863385, -6.35848597656913, -2.059108837045416),
    @(2600, "walkingToRunning", 14.78357696533203, -19.61955070495605, 58.06728363037109, -9.929815484465564, -8.435556281384915, 4.079216970814213),
    @(2700, "walkingToRunning", -78.68696594238281, -17.28630638122559, -26.55977630615234, -5.980313006064948, -16.70210841062257, 0.5289435215133755),
    @(2800, "walkingToRunning", -3.455702781677246, -1.554847240447998, -16.96383094787598, 5.703074371214343, -13.00758359758107, -1.666179502181835),
    @(2900, "walkingToRunning", -10.43552684783936, -8.621312141418457, -5.604441642761231, 3.259694958762311, -1.829939403122192, -3.400759592759585)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}